# Replace the "nihk" label cell with the new array-style field placeholder
# "{base.name}" and point the selection at the header row of the
# template-field table (A4) instead of the previously selected C8.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A1 used to hold the shared string "nihk" - swap it for the new
# placeholder string used by the array/json-to-hashmap write path.
$ws.Range("A1").Value = "{base.name}"

# Move the active selection to A4 (first cell of the field-name row).
$ws.Range("A4").Select()
